$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = "TEST-28845"
$ws.Range("B3").Value = "100;rated-6009;1"
$ws.Range("I3").Value = "P1;100;P1"
$ws.Range("P3").Value = "2s"
$ws.Range("Q3").Value = "4d"
$ws.Range("R3").Value = "3s"
$ws.Range("S3").Value = "4d"
$ws.Range("V3").Value = "P1"

# Row 4
$ws.Range("A4").Value = "TEST-14007"
$ws.Range("B4").Value = "1000-2;anon"
$ws.Range("I4").Value = "P1;1000;B3"
$ws.Range("J4").Value = "P1;1000;P2"
$ws.Range("Q4").Value = "4s"
$ws.Range("P4").Value = "Ah"
$ws.Range("R4").Value = "2s"
$ws.Range("S4").Value = "Ah"
$ws.Range("V4").Value = "B3;P2"

# Update view state to match target
$ws.Application.ActiveWindow.ScrollColumn = 11
$ws.Range("V7").Select()
